$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "nuanced details across their entire",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "subtle details across the entire",
    2
)
